$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data rows 2 and 3
$ws.Range("B2").Value = 246.57255000000001
$ws.Range("C2").Value = 0.16050700000000001
$ws.Range("D2").Value = 3188

$ws.Range("B3").Value = 477.22721200000001
$ws.Range("C3").Value = 0.68484199999999995
$ws.Range("D3").Value = 2447

# Fill previously-empty row 4 with new data (new shared string "bs3")
$ws.Range("A4").Value = "bs3"
$ws.Range("B4").Value = 579.68972900000006
$ws.Range("C4").Value = 0.16483600000000001
$ws.Range("D4").Value = 3101

# Update the active selection shown in the sheet view
$ws.Range("G6").Select()
